# Modified format tables in: Ambits_sve, Ambits_rang, Ambits_Rang_SVE_mes, Ambits_verificació
#
# The historical table on the active sheet gains a new "SE3" bucket:
#   - two new trailing columns, AJ ("SE3_CI") and AK ("SE3_CNI"), with their
#     per-region totals;
#   - the previous trailing columns, AH ("SE2_CI") and AI ("SE2_CNI"), switch
#     from numeric storage to text storage (same literal "<n>.0" figures),
#     matching how the refreshed export represents them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers -----------------------------------------------
$ws.Range("AJ1").Value = "SE3_CI"
$ws.Range("AK1").Value = "SE3_CNI"

# --- Data rows 2..11 ----------------------------------------------------
# AH/AI keep the exact same figures they already hold, just re-typed as text.
$ahValues = @("2779.0","2921.0","3578.0","1326.0","1770.0","863.0","1375.0","511.0","797.0","15920.0")
$aiValues = @("3381.0","2865.0","2746.0","779.0","736.0","562.0","1021.0","435.0","729.0","13254.0")

# New SE3_CI / SE3_CNI figures per region (row 2..11, last row is the Totals row).
$ajValues = @(2636.0,2231.0,2382.0,841.0,1322.0,538.0,1100.0,432.0,522.0,12004.0)
$akValues = @(2586.0,1726.0,2312.0,404.0,652.0,434.0,679.0,287.0,537.0,9617.0)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2

    # Write AH/AI as formulas producing the literal text, then collapse them
    # to plain text values below (keeps default "General" cell styling).
    $ws.Cells.Item($row, 34).Formula = '="' + $ahValues[$i] + '"'
    $ws.Cells.Item($row, 35).Formula = '="' + $aiValues[$i] + '"'

    $ws.Cells.Item($row, 36).Value = $ajValues[$i]
    $ws.Cells.Item($row, 37).Value = $akValues[$i]
}

# Collapse the AH:AI formulas into literal text values (no residual formulas,
# no number re-inference, no new number-format style).
$rng = $ws.Range("AH2:AI11")
$rng.Copy()
$rng.PasteSpecial(-4163)
